$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @(5,5,5,5,-1,0,5,-1,5,5,-1,-1,5,5,0,0,4,3,5,5,0,2,1,0,-1,0,0,5,0,5,5,0,5,0,5,1,5,5,-1,5,1,-1,5,0,5,-1,5,-1,4,-1,0,-1,-1,0,2,-1,-1,2,5,5,0,5,-1,4,5,5,5,5,-1,-1,-1,5,5,5,5,5,4,5,5,5,5,5,5,5,0,4,5,5,4,5,5,0,0,5,0,4,5,5,3,-1,-1,5,5,2,0,0,5,5,3,2,2,5,-1,0,5,-1,0,0,0,0,-1,5,5,2,5,5,0,2,2,5,5,5,0,5,5,5,5,5,5,5,-1,-1,5,0,-1,2,5,5,5,3,2,4,2,2,5,5,2,2,4,5,-1,2,2,2,-1,0,-1,5,4,0,-1,4,5,3,4,3,0,2,-1)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $labels[$i]
}
